$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J12").Value = 0.2150495036779461
$ws.Range("I13").Value = 0.24
$ws.Range("H14").Value = 0.3087982760018804
$ws.Range("G15").Value = 0.32
$ws.Range("F16").Value = 0.4476495795507702
$ws.Range("E17").Value = 0.1088966743764388
$ws.Range("D18").Value = 0.1461563307127136
$ws.Range("C19").Value = 0.09547648014918764
$ws.Range("B20").Value = 0.0959495356205764
